# Fruta / hortaliza, semanal
# Insert a new week's worth of data (two rows: "Primera" and "Segunda" quality)
# right above the current row 400, pushing the existing rows 400-427 down to
# 402-429.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 400 (existing row 400 and everything
# below shifts down by 2).
$ws.Rows.Item(400).Insert()
$ws.Rows.Item(400).Insert()

# New row 400: "Primera" quality entry for the latest week.
$ws.Range("A400").Value = 11
$ws.Range("B400").Value = "Vega Monumental Concepción"
$ws.Range("C400").Value = "Bíobío"
$ws.Range("D400").Value = 45106
$ws.Range("E400").Value = 8
$ws.Range("F400").Value = 100112009
$ws.Range("G400").Value = "Acelga"
$ws.Range("H400").Value = "Sin especificar"
$ws.Range("I400").Value = "Primera"
$ws.Range("J400").Value = 200
$ws.Range("K400").Value = 600
$ws.Range("L400").Value = 700
$ws.Range("M400").Value = 650
$ws.Range("N400").Value = "`$/atado"
$ws.Range("O400").Value = "Región de Ñuble"
$ws.Range("P400").Value = 650
$ws.Range("Q400").Value = 1
$ws.Range("R400").Value = "Hortaliza"

# New row 401: "Segunda" quality entry for the latest week.
$ws.Range("A401").Value = 11
$ws.Range("B401").Value = "Vega Monumental Concepción"
$ws.Range("C401").Value = "Bíobío"
$ws.Range("D401").Value = 45106
$ws.Range("E401").Value = 8
$ws.Range("F401").Value = 100112009
$ws.Range("G401").Value = "Acelga"
$ws.Range("H401").Value = "Sin especificar"
$ws.Range("I401").Value = "Segunda"
$ws.Range("J401").Value = 100
$ws.Range("K401").Value = 500
$ws.Range("L401").Value = 500
$ws.Range("M401").Value = 500
$ws.Range("N401").Value = "`$/atado"
$ws.Range("O401").Value = "Región de Ñuble"
$ws.Range("P401").Value = 500
$ws.Range("Q401").Value = 1
$ws.Range("R401").Value = "Hortaliza"
